$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 59, shifting existing rows 59-85
# down to become rows 61-87.
$ws.Range("A59:A60").EntireRow.Insert()

# Fill in new row 59
$ws.Range("A59").Value = 1
$ws.Range("B59").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C59").Value = "Arica y Parinacota"
$ws.Range("D59").Value = 44992
$ws.Range("E59").Value = 15
$ws.Range("F59").Value = "Fruta"
$ws.Range("G59").Value = 100103
$ws.Range("H59").Value = "Frutos de hueso (carozo)"
$ws.Range("I59").Value = 100103004
$ws.Range("J59").Value = "Durazno"
$ws.Range("K59").Value = "Phillips Cling"
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 250
$ws.Range("N59").Value = 20000
$ws.Range("O59").Value = 22000
$ws.Range("P59").Value = 20800
$ws.Range("Q59").Value = "$/bandeja 18 kilos granel"
$ws.Range("R59").Value = "Región de O'Higgins"
$ws.Range("S59").Value = 1156
$ws.Range("T59").Value = 18

# Fill in new row 60
$ws.Range("A60").Value = 1
$ws.Range("B60").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C60").Value = "Arica y Parinacota"
$ws.Range("D60").Value = 44992
$ws.Range("E60").Value = 15
$ws.Range("F60").Value = "Fruta"
$ws.Range("G60").Value = 100103
$ws.Range("H60").Value = "Frutos de hueso (carozo)"
$ws.Range("I60").Value = 100103004
$ws.Range("J60").Value = "Durazno"
$ws.Range("K60").Value = "September Sweet"
$ws.Range("L60").Value = "Primera"
$ws.Range("M60").Value = 250
$ws.Range("N60").Value = 20000
$ws.Range("O60").Value = 22000
$ws.Range("P60").Value = 20800
$ws.Range("Q60").Value = "$/bandeja 18 kilos granel"
$ws.Range("R60").Value = "Región de O'Higgins"
$ws.Range("S60").Value = 1156
$ws.Range("T60").Value = 18
